$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (GHA / E / Final / PCM / eta.fu) was missing a value in column H
# (H1971). Fill it in with the same value already present in column I
# (I2000), matching the format/style already used by the rest of row 4.
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = 0.20099155835454899

# Move the active selection from H5 to the newly completed H4 cell.
$ws.Range("H4").Select()
